# Adding 2 search test cases to the "Test Cases" sheet (rows 79 and 80).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 79: TestCase_B78 / OPQA-412 -----------------------------------
$ws.Range("A79").Value = "TestCase_B78"
$ws.Range("B79").Value = "OPQA-412"
$ws.Range("C79").Value = "Verify that nothing gets displayed in the search type ahead if search query is not interpreted by the system"
$ws.Range("D79").Value = "Y"
$ws.Range("E79").Value = "SKIP"

# --- New row 80: TestCase_B79 / OPQA-392 -----------------------------------
$ws.Range("A80").Value = "TestCase_B79"
$ws.Range("B80").Value = "OPQA-392"
$ws.Range("C80").Value = "Verify that the following changes take place when user clicks on any CATEGORIES option in the search type ahead while ALL option is selected in the search drop down:`na)Correct keyword gets displayed in the search box`nb)ARTICLES option gets selected both in the search drop down and left navigation pane`nc)Only articles get displayed in the search results page`nd)Correct category gets selected in CATEGORIES filter in the left navigation pane with filter in expanded state"
$ws.Range("D80").Value = "Y"
$ws.Range("E80").Value = "PASS"

# Row 80 description is long -> wraps onto many lines, same as the other
# multi-step rows in this sheet (e.g. row 32 which also uses ht="90").
$ws.Range("C80").WrapText = $true
$ws.Rows.Item(80).RowHeight = 90

# Row 74-78 ("ALL option" test cases) were marked PASS but should now read
# SKIP, matching the rest of the suite.
$ws.Range("E74").Value = "SKIP"
$ws.Range("E75").Value = "SKIP"
$ws.Range("E76").Value = "SKIP"
$ws.Range("E77").Value = "SKIP"
$ws.Range("E78").Value = "SKIP"

# Keep the visible selection/viewport in sync with the newly added rows.
$ws.Range("D2:D80").Select()
